$d = $word.ActiveDocument

# "Good " + "afternoon" -> "This is A1 batch." (merges the two runs into one)
$d.Content.Find.Execute("Good afternoon", $false, $false, $false, $false, $false, $true, 1, $false, "This is A1 batch.", 2)

# Drop the trailing "…!" run entirely
$d.Content.Find.Execute("…!", $false, $false, $false, $false, $false, $true, 1, $false, "", 2)
